# feat: add 2022-Q3 data
#
# 1) "总计" (summary) sheet: insert a new top data row for 2022-Q3 and
#    shift the previously-existing quarters down by one row.
# 2) Insert a brand-new worksheet "2022-Q3" (fund holdings detail) right
#    before the existing "2022-Q2" sheet, with the same layout/format as
#    the other quarterly sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: update "总计" summary sheet (first sheet)
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summaryData = @(
    @("2022-Q3", 16, 2.23),
    @("2022-Q2", 12, 2.04),
    @("2022-Q1", 20, 2.08),
    @("2021-Q4", 23, 5.03),
    @("2021-Q3", 16, 4.88),
    @("2021-Q2", 17, 8.62),
    @("2021-Q1", 26, 14.11),
    @("2020-Q4", 8, 1.68)
)

$r = 2
foreach ($row in $summaryData) {
    $summary.Cells.Item($r, 1).Value = $r - 2
    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = $row[2]
    $r++
}

# The table grew by one row (A1:D8 -> A1:D9); make sure the newly
# created last row's index cell (column A) carries the same format as
# the rest of the index column.
$summary.Range("A2").Copy() | Out-Null
$summary.Range("A9").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# Step 2: insert the new "2022-Q3" worksheet before "2022-Q2"
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($anchor)
$q3.Name = "2022-Q3"

# re-fetch the old "2022-Q2" sheet fresh (it shifted to index 3, and the
# previous $anchor handle now tracks the newly-inserted sheet instead)
$template = $wb.Worksheets.Item(3)

# Match sheet-level outline flags used by every other sheet in this
# workbook.
$q3.Outline.SummaryRow = 1
$q3.Outline.SummaryColumn = 1

# Header row (B1:H1), copying the bold/border/center format used on
# every other quarterly sheet.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$template.Range("B1:H1").Copy() | Out-Null
$q3.Range("B1:H1").PasteSpecial(-4122) | Out-Null

$rows = @(
    @("501311", "嘉实恒生港股通新经济指数（LOF）A", "15.36", "94.00", "4.60", "0.7066", 6),
    @("006614", "嘉实恒生港股通新经济指数（LOF）C", "14.03", "94.00", "4.60", "0.6454", 6),
    @("009007", "兴全沪港深两年持有期混合", "16.31", "92.09", "3.49", "0.5692", 9),
    @("011924", "嘉实港股互联网产业核心资产混合A", "1.20", "87.88", "3.84", "0.0461", 9),
    @("513960", "博时港股通消费ETF", "0.79", "97.92", "4.02", "0.0318", 7),
    @("006787", "泰康中证港股通大消费主题指数C", "0.80", "87.17", "3.92", "0.0314", 6),
    @("159735", "银华中证港股通消费主题ETF", "0.75", "92.71", "3.80", "0.0285", 7),
    @("007151", "前海开源沪港深聚瑞混合", "0.53", "82.80", "5.07", "0.0269", 10),
    @("513320", "易方达恒生港股通新经济ETF", "0.49", "98.67", "4.86", "0.0238", 6),
    @("513230", "华夏中证港股通消费主题ETF", "0.50", "96.82", "4.00", "0.0200", 7),
    @("513070", "易方达中证港股通消费主题ETF", "0.49", "96.86", "3.95", "0.0194", 7),
    @("006786", "泰康中证港股通大消费主题指数A", "0.45", "87.17", "3.92", "0.0176", 6),
    @("517880", "华泰柏瑞中证沪港深品牌消费50ETF", "0.47", "92.73", "3.73", "0.0175", 8),
    @("513590", "鹏华中证港股通消费主题ETF", "0.52", "79.92", "3.31", "0.0172", 7),
    @("011925", "嘉实港股互联网产业核心资产混合C", "0.41", "87.88", "3.84", "0.0157", 9),
    @("162416", "华宝港股通恒生香港35指数（LOF）", "0.20", "93.77", "4.14", "0.0083", 9)
)

$lastRow = 1 + $rows.Length

# Columns B..G hold text that looks numeric (fund codes with leading
# zeros, percentages formatted to two decimals, etc.) - force text so
# Excel doesn't silently coerce them to numbers.
$q3.Range("B2:G$lastRow").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $r - 2
    $q3.Cells.Item($r, 2).Value = $row[0]
    $q3.Cells.Item($r, 3).Value = $row[1]
    $q3.Cells.Item($r, 4).Value = $row[2]
    $q3.Cells.Item($r, 5).Value = $row[3]
    $q3.Cells.Item($r, 6).Value = $row[4]
    $q3.Cells.Item($r, 7).Value = $row[5]
    $q3.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# Column A (row index 0..15) uses the same bold/border/center format as
# every other sheet's index column.
$template.Range("A2").Copy() | Out-Null
$q3.Range("A2:A$lastRow").PasteSpecial(-4122) | Out-Null

Write-Output "2022-Q3 sheet added and 总计 refreshed"
